$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns E, F, G (whole columns) - this shifts H->E and I->F,
# reusing the existing shared-string values and cell styles that were in H/I.
$ws.Range("E1:G1").EntireColumn.Delete()

# Fill in column D (the "who" column) which used to be blank for several rows.
$ws.Range("D2").Value = "Хозяин"
$ws.Range("D5").Value = "Михаил"
$ws.Range("D6").Value = "Марина"
$ws.Range("D7").Value = "Илья"

# Row 4 did not previously have a gender cell (old I4); add the new one.
$ws.Range("F4").Value = "M"

# Update the current selection to match the edited workbook.
$ws.Range("M12").Select()

Write-Output "done"
